# Categories.xlsx refactor: "refactored parsing with multiple strategies"
#
# Semantic summary of the change (derived from the OOXML diff):
#  - "Business" (the COMMON category) was renamed to "Economics"; all rows
#    that used "Business" as their parent now use "Economics".
#  - "Small Business" was renamed to plain "Business" (a *new* COMMON
#    top-level category, distinct from the old "Business"/now "Economics").
#  - "Arts" was renamed to "Culture"; all rows that used "Arts" as their
#    parent now use "Culture".
#  - 14 new category rows were appended (rows 71-84): Theatre, Science,
#    History (children of Culture/COMMON), eight European countries
#    (children of Europe/REGION), and Animals (child of Ecology/COMMON).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Business" -> "Economics" (row 29 is the category itself;
#     rows 31-37,40-42 are children that pointed at it as parent) ---
$ws.Cells.Item(29,1).Value = "Economics"
$ws.Cells.Item(31,3).Value = "Economics"
$ws.Cells.Item(32,3).Value = "Economics"
$ws.Cells.Item(33,3).Value = "Economics"
$ws.Cells.Item(34,3).Value = "Economics"
$ws.Cells.Item(35,3).Value = "Economics"
$ws.Cells.Item(36,3).Value = "Economics"
$ws.Cells.Item(37,3).Value = "Economics"
$ws.Cells.Item(40,3).Value = "Economics"
$ws.Cells.Item(41,3).Value = "Economics"
$ws.Cells.Item(42,3).Value = "Economics"

# --- Rename "Small Business" -> "Business" (row 42 is the category itself) ---
$ws.Cells.Item(42,1).Value = "Business"

# --- Rename "Arts" -> "Culture" (row 51 is the category itself; rows
#     52-58 and 70 are children that pointed at it as parent) ---
$ws.Cells.Item(51,1).Value = "Culture"
$ws.Cells.Item(52,3).Value = "Culture"
$ws.Cells.Item(53,3).Value = "Culture"
$ws.Cells.Item(54,3).Value = "Culture"
$ws.Cells.Item(55,3).Value = "Culture"
$ws.Cells.Item(56,3).Value = "Culture"
$ws.Cells.Item(57,3).Value = "Culture"
$ws.Cells.Item(58,3).Value = "Culture"
$ws.Cells.Item(70,3).Value = "Culture"

# --- Append 14 new rows (71-84). First clone the formatting of existing
#     data rows (style index 1) onto the new block so the new cells match
#     the sheet's existing look. Row 72 only has a name+type (no parent),
#     like rows such as 22-24, so its format is copied from a same-shaped
#     two-column row (A22:B22) to avoid materialising a spurious empty
#     C72 cell; the rest are copied from a fully-populated three-column
#     row (A4:C4). ---
$ws.Range("A4:C4").Copy()
$ws.Range("A71:C71").PasteSpecial(-4122)
$ws.Range("A73:C84").PasteSpecial(-4122)
$ws.Range("A22:B22").Copy()
$ws.Range("A72:B72").PasteSpecial(-4122)

$ws.Cells.Item(71,1).Value = "Theatre"
$ws.Cells.Item(71,2).Value = "COMMON"
$ws.Cells.Item(71,3).Value = "Culture"

$ws.Cells.Item(72,1).Value = "Science"
$ws.Cells.Item(72,2).Value = "COMMON"

$ws.Cells.Item(73,1).Value = "History"
$ws.Cells.Item(73,2).Value = "COMMON"
$ws.Cells.Item(73,3).Value = "Culture"

$ws.Cells.Item(74,1).Value = "Italy"
$ws.Cells.Item(74,2).Value = "REGION"
$ws.Cells.Item(74,3).Value = "Europe"

$ws.Cells.Item(75,1).Value = "France"
$ws.Cells.Item(75,2).Value = "REGION"
$ws.Cells.Item(75,3).Value = "Europe"

$ws.Cells.Item(76,1).Value = "Greece"
$ws.Cells.Item(76,2).Value = "REGION"
$ws.Cells.Item(76,3).Value = "Europe"

$ws.Cells.Item(77,1).Value = "Germany"
$ws.Cells.Item(77,2).Value = "REGION"
$ws.Cells.Item(77,3).Value = "Europe"

$ws.Cells.Item(78,1).Value = "Denmark"
$ws.Cells.Item(78,2).Value = "REGION"
$ws.Cells.Item(78,3).Value = "Europe"

$ws.Cells.Item(79,1).Value = "Sweden"
$ws.Cells.Item(79,2).Value = "REGION"
$ws.Cells.Item(79,3).Value = "Europe"

$ws.Cells.Item(80,1).Value = "Norway"
$ws.Cells.Item(80,2).Value = "REGION"
$ws.Cells.Item(80,3).Value = "Europe"

$ws.Cells.Item(81,1).Value = "Finland"
$ws.Cells.Item(81,2).Value = "REGION"
$ws.Cells.Item(81,3).Value = "Europe"

$ws.Cells.Item(82,1).Value = "Spain"
$ws.Cells.Item(82,2).Value = "REGION"
$ws.Cells.Item(82,3).Value = "Europe"

$ws.Cells.Item(83,1).Value = "Austria"
$ws.Cells.Item(83,2).Value = "REGION"
$ws.Cells.Item(83,3).Value = "Europe"

$ws.Cells.Item(84,1).Value = "Animals"
$ws.Cells.Item(84,2).Value = "COMMON"
$ws.Cells.Item(84,3).Value = "Ecology"
